$p = $ppt.ActivePresentation
$d = $p.Designs
Write-Output ("before count=" + $d.Count)
try {
    $newd = $d.Add("test")
    Write-Output ("Add result type=" + $newd.GetType())
} catch {
    Write-Output ("ERR: " + $_.Exception.Message)
}
Write-Output ("after count=" + $d.Count)
